# Commit: "Fruta / hortaliza, semanal" -- weekly refresh of the Lechuga price sheet.
# A new week (2022-01-17, serial 44578) of "Escarola" Primera/Segunda quotes is inserted
# at the top of the date-ordered data block (old rows 568-608 shift down by 2 to 570-610),
# and the two oldest rows that fall out of the tracked window are re-appended at the very
# end of the sheet (rows 609-610, identical to what used to be rows 607-608).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing data block down by 2 rows (write the new values in place) ---
# Row 568
$ws.Cells.Item(568, 4).Value = 44578

# Row 569
$ws.Cells.Item(569, 4).Value = 44578

# Row 570
$ws.Cells.Item(570, 4).Value = 44490
$ws.Cells.Item(570, 11).Value = 5000
$ws.Cells.Item(570, 12).Value = 6000
$ws.Cells.Item(570, 13).Value = 5500
$ws.Cells.Item(570, 16).Value = 458

# Row 571
$ws.Cells.Item(571, 4).Value = 44490
$ws.Cells.Item(571, 11).Value = 5000
$ws.Cells.Item(571, 12).Value = 6000
$ws.Cells.Item(571, 13).Value = 5500
$ws.Cells.Item(571, 16).Value = 306

# Row 572
$ws.Cells.Item(572, 9).Value = 'Primera'
$ws.Cells.Item(572, 11).Value = 11000
$ws.Cells.Item(572, 12).Value = 12000
$ws.Cells.Item(572, 13).Value = 11500
$ws.Cells.Item(572, 14).Value = '$/caja 12 unidades'
$ws.Cells.Item(572, 16).Value = 958
$ws.Cells.Item(572, 17).Value = 12

# Row 573
$ws.Cells.Item(573, 8).Value = 'Escarola'
$ws.Cells.Item(573, 9).Value = 'Segunda'
$ws.Cells.Item(573, 10).Value = 120
$ws.Cells.Item(573, 11).Value = 11000
$ws.Cells.Item(573, 12).Value = 12000
$ws.Cells.Item(573, 13).Value = 11500
$ws.Cells.Item(573, 14).Value = '$/caja 18 unidades'
$ws.Cells.Item(573, 16).Value = 639
$ws.Cells.Item(573, 17).Value = 18

# Row 574
$ws.Cells.Item(574, 8).Value = 'Escarola'
$ws.Cells.Item(574, 9).Value = 'Tercera'
$ws.Cells.Item(574, 14).Value = '$/caja 24 unidades'
$ws.Cells.Item(574, 16).Value = 365
$ws.Cells.Item(574, 17).Value = 24

# Row 575
$ws.Cells.Item(575, 4).Value = 44427
$ws.Cells.Item(575, 8).Value = 'Marina'
$ws.Cells.Item(575, 10).Value = 140
$ws.Cells.Item(575, 11).Value = 8500
$ws.Cells.Item(575, 12).Value = 9000
$ws.Cells.Item(575, 13).Value = 8750
$ws.Cells.Item(575, 16).Value = 729

# Row 576
$ws.Cells.Item(576, 4).Value = 44427
$ws.Cells.Item(576, 8).Value = 'Marina'
$ws.Cells.Item(576, 11).Value = 8500
$ws.Cells.Item(576, 12).Value = 9000
$ws.Cells.Item(576, 13).Value = 8750
$ws.Cells.Item(576, 16).Value = 486

# Row 577
$ws.Cells.Item(577, 8).Value = 'Escarola'
$ws.Cells.Item(577, 11).Value = 6000
$ws.Cells.Item(577, 12).Value = 7000
$ws.Cells.Item(577, 13).Value = 6500
$ws.Cells.Item(577, 16).Value = 542

# Row 578
$ws.Cells.Item(578, 8).Value = 'Escarola'
$ws.Cells.Item(578, 10).Value = 120
$ws.Cells.Item(578, 11).Value = 6000
$ws.Cells.Item(578, 12).Value = 7000
$ws.Cells.Item(578, 13).Value = 6500
$ws.Cells.Item(578, 16).Value = 361

# Row 579
$ws.Cells.Item(579, 4).Value = 44491
$ws.Cells.Item(579, 8).Value = 'Marina'
$ws.Cells.Item(579, 11).Value = 4000
$ws.Cells.Item(579, 13).Value = 4500
$ws.Cells.Item(579, 16).Value = 375

# Row 580
$ws.Cells.Item(580, 4).Value = 44491
$ws.Cells.Item(580, 8).Value = 'Marina'
$ws.Cells.Item(580, 10).Value = 140
$ws.Cells.Item(580, 11).Value = 4000
$ws.Cells.Item(580, 13).Value = 4500
$ws.Cells.Item(580, 16).Value = 250

# Row 581
$ws.Cells.Item(581, 4).Value = 44293
$ws.Cells.Item(581, 8).Value = 'Escarola'
$ws.Cells.Item(581, 10).Value = 120
$ws.Cells.Item(581, 11).Value = 4500
$ws.Cells.Item(581, 12).Value = 5000
$ws.Cells.Item(581, 13).Value = 4750
$ws.Cells.Item(581, 14).Value = '$/caja 12 unidades'
$ws.Cells.Item(581, 16).Value = 396
$ws.Cells.Item(581, 17).Value = 12

# Row 582
$ws.Cells.Item(582, 4).Value = 44293
$ws.Cells.Item(582, 8).Value = 'Escarola'
$ws.Cells.Item(582, 9).Value = 'Segunda'
$ws.Cells.Item(582, 11).Value = 4500
$ws.Cells.Item(582, 12).Value = 5000
$ws.Cells.Item(582, 13).Value = 4750
$ws.Cells.Item(582, 14).Value = '$/caja 18 unidades'
$ws.Cells.Item(582, 16).Value = 264
$ws.Cells.Item(582, 17).Value = 18

# Row 583
$ws.Cells.Item(583, 4).Value = 44266
$ws.Cells.Item(583, 8).Value = 'Conconina(o)'
$ws.Cells.Item(583, 9).Value = 'Primera'
$ws.Cells.Item(583, 10).Value = 1250
$ws.Cells.Item(583, 11).Value = 800
$ws.Cells.Item(583, 12).Value = 900
$ws.Cells.Item(583, 13).Value = 850
$ws.Cells.Item(583, 14).Value = '$/unidad'
$ws.Cells.Item(583, 16).Value = 850
$ws.Cells.Item(583, 17).Value = 1

# Row 584
$ws.Cells.Item(584, 4).Value = 44264
$ws.Cells.Item(584, 8).Value = 'Marina'
$ws.Cells.Item(584, 11).Value = 10000
$ws.Cells.Item(584, 12).Value = 11000
$ws.Cells.Item(584, 13).Value = 10500
$ws.Cells.Item(584, 16).Value = 875

# Row 585
$ws.Cells.Item(585, 4).Value = 44264
$ws.Cells.Item(585, 8).Value = 'Marina'
$ws.Cells.Item(585, 11).Value = 10000
$ws.Cells.Item(585, 12).Value = 11000
$ws.Cells.Item(585, 13).Value = 10500
$ws.Cells.Item(585, 16).Value = 583

# Row 586
$ws.Cells.Item(586, 4).Value = 44494
$ws.Cells.Item(586, 8).Value = 'Escarola'
$ws.Cells.Item(586, 10).Value = 120
$ws.Cells.Item(586, 11).Value = 6000
$ws.Cells.Item(586, 12).Value = 7000
$ws.Cells.Item(586, 13).Value = 6500
$ws.Cells.Item(586, 14).Value = '$/caja 12 unidades'
$ws.Cells.Item(586, 16).Value = 542
$ws.Cells.Item(586, 17).Value = 12

# Row 587
$ws.Cells.Item(587, 4).Value = 44494
$ws.Cells.Item(587, 8).Value = 'Escarola'
$ws.Cells.Item(587, 10).Value = 120
$ws.Cells.Item(587, 11).Value = 6000
$ws.Cells.Item(587, 12).Value = 7000
$ws.Cells.Item(587, 13).Value = 6500
$ws.Cells.Item(587, 14).Value = '$/caja 18 unidades'
$ws.Cells.Item(587, 16).Value = 361
$ws.Cells.Item(587, 17).Value = 18

# Row 588
$ws.Cells.Item(588, 8).Value = 'Conconina(o)'
$ws.Cells.Item(588, 10).Value = 700
$ws.Cells.Item(588, 11).Value = 700
$ws.Cells.Item(588, 12).Value = 800
$ws.Cells.Item(588, 13).Value = 750
$ws.Cells.Item(588, 14).Value = '$/unidad'
$ws.Cells.Item(588, 16).Value = 750
$ws.Cells.Item(588, 17).Value = 1

# Row 589
$ws.Cells.Item(589, 8).Value = 'Conconina(o)'
$ws.Cells.Item(589, 10).Value = 1000
$ws.Cells.Item(589, 11).Value = 500
$ws.Cells.Item(589, 12).Value = 600
$ws.Cells.Item(589, 13).Value = 550
$ws.Cells.Item(589, 14).Value = '$/unidad'
$ws.Cells.Item(589, 16).Value = 550
$ws.Cells.Item(589, 17).Value = 1

# Row 590
$ws.Cells.Item(590, 4).Value = 44390
$ws.Cells.Item(590, 10).Value = 150
$ws.Cells.Item(590, 11).Value = 5500
$ws.Cells.Item(590, 12).Value = 6000
$ws.Cells.Item(590, 13).Value = 5750
$ws.Cells.Item(590, 16).Value = 479

# Row 591
$ws.Cells.Item(591, 4).Value = 44390
$ws.Cells.Item(591, 10).Value = 136
$ws.Cells.Item(591, 11).Value = 5500
$ws.Cells.Item(591, 12).Value = 6000
$ws.Cells.Item(591, 13).Value = 5750
$ws.Cells.Item(591, 16).Value = 319

# Row 592
$ws.Cells.Item(592, 8).Value = 'Escarola'
$ws.Cells.Item(592, 11).Value = 7500
$ws.Cells.Item(592, 12).Value = 8000
$ws.Cells.Item(592, 13).Value = 7750
$ws.Cells.Item(592, 16).Value = 646

# Row 593
$ws.Cells.Item(593, 8).Value = 'Escarola'
$ws.Cells.Item(593, 11).Value = 7500
$ws.Cells.Item(593, 12).Value = 8000
$ws.Cells.Item(593, 13).Value = 7750
$ws.Cells.Item(593, 16).Value = 431

# Row 594
$ws.Cells.Item(594, 4).Value = 44279
$ws.Cells.Item(594, 8).Value = 'Marina'
$ws.Cells.Item(594, 10).Value = 120
$ws.Cells.Item(594, 11).Value = 10000
$ws.Cells.Item(594, 12).Value = 11000
$ws.Cells.Item(594, 13).Value = 10500
$ws.Cells.Item(594, 16).Value = 875

# Row 595
$ws.Cells.Item(595, 4).Value = 44279
$ws.Cells.Item(595, 8).Value = 'Marina'
$ws.Cells.Item(595, 11).Value = 10000
$ws.Cells.Item(595, 12).Value = 11000
$ws.Cells.Item(595, 13).Value = 10500
$ws.Cells.Item(595, 16).Value = 583

# Row 596
$ws.Cells.Item(596, 4).Value = 44481
$ws.Cells.Item(596, 10).Value = 130
$ws.Cells.Item(596, 11).Value = 6000
$ws.Cells.Item(596, 12).Value = 7000
$ws.Cells.Item(596, 13).Value = 6500
$ws.Cells.Item(596, 16).Value = 542

# Row 597
$ws.Cells.Item(597, 4).Value = 44481
$ws.Cells.Item(597, 11).Value = 6000
$ws.Cells.Item(597, 12).Value = 7000
$ws.Cells.Item(597, 13).Value = 6500
$ws.Cells.Item(597, 16).Value = 361

# Row 598
$ws.Cells.Item(598, 4).Value = 44277
$ws.Cells.Item(598, 8).Value = 'Escarola'
$ws.Cells.Item(598, 9).Value = 'Primera'
$ws.Cells.Item(598, 10).Value = 120
$ws.Cells.Item(598, 11).Value = 7000
$ws.Cells.Item(598, 12).Value = 8000
$ws.Cells.Item(598, 13).Value = 7500
$ws.Cells.Item(598, 14).Value = '$/caja 12 unidades'
$ws.Cells.Item(598, 16).Value = 625
$ws.Cells.Item(598, 17).Value = 12

# Row 599
$ws.Cells.Item(599, 4).Value = 44277
$ws.Cells.Item(599, 9).Value = 'Segunda'
$ws.Cells.Item(599, 10).Value = 120
$ws.Cells.Item(599, 11).Value = 7000
$ws.Cells.Item(599, 12).Value = 8000
$ws.Cells.Item(599, 13).Value = 7500
$ws.Cells.Item(599, 14).Value = '$/caja 18 unidades'
$ws.Cells.Item(599, 16).Value = 417
$ws.Cells.Item(599, 17).Value = 18

# Row 600
$ws.Cells.Item(600, 8).Value = 'Conconina(o)'
$ws.Cells.Item(600, 10).Value = 1500
$ws.Cells.Item(600, 11).Value = 450
$ws.Cells.Item(600, 12).Value = 500
$ws.Cells.Item(600, 13).Value = 475
$ws.Cells.Item(600, 14).Value = '$/unidad'
$ws.Cells.Item(600, 16).Value = 475
$ws.Cells.Item(600, 17).Value = 1

# Row 601
$ws.Cells.Item(601, 4).Value = 44525
$ws.Cells.Item(601, 10).Value = 140
$ws.Cells.Item(601, 11).Value = 2500
$ws.Cells.Item(601, 12).Value = 3000
$ws.Cells.Item(601, 13).Value = 2750
$ws.Cells.Item(601, 16).Value = 229

# Row 602
$ws.Cells.Item(602, 4).Value = 44525
$ws.Cells.Item(602, 10).Value = 160
$ws.Cells.Item(602, 11).Value = 2500
$ws.Cells.Item(602, 12).Value = 3000
$ws.Cells.Item(602, 13).Value = 2750
$ws.Cells.Item(602, 16).Value = 153

# Row 603
$ws.Cells.Item(603, 8).Value = 'Escarola'
$ws.Cells.Item(603, 11).Value = 5000
$ws.Cells.Item(603, 12).Value = 6000
$ws.Cells.Item(603, 13).Value = 5500
$ws.Cells.Item(603, 16).Value = 458

# Row 604
$ws.Cells.Item(604, 8).Value = 'Escarola'
$ws.Cells.Item(604, 11).Value = 5000
$ws.Cells.Item(604, 12).Value = 6000
$ws.Cells.Item(604, 13).Value = 5500
$ws.Cells.Item(604, 16).Value = 306

# Row 605
$ws.Cells.Item(605, 4).Value = 44327
$ws.Cells.Item(605, 8).Value = 'Marina'
$ws.Cells.Item(605, 10).Value = 120
$ws.Cells.Item(605, 11).Value = 4000
$ws.Cells.Item(605, 12).Value = 4500
$ws.Cells.Item(605, 13).Value = 4250
$ws.Cells.Item(605, 16).Value = 354

# Row 606
$ws.Cells.Item(606, 4).Value = 44327
$ws.Cells.Item(606, 8).Value = 'Marina'
$ws.Cells.Item(606, 10).Value = 120
$ws.Cells.Item(606, 11).Value = 4000
$ws.Cells.Item(606, 12).Value = 4500
$ws.Cells.Item(606, 13).Value = 4250
$ws.Cells.Item(606, 16).Value = 236

# Row 607
$ws.Cells.Item(607, 4).Value = 44384
$ws.Cells.Item(607, 10).Value = 160

# Row 608
$ws.Cells.Item(608, 4).Value = 44384
$ws.Cells.Item(608, 10).Value = 200
$ws.Cells.Item(608, 13).Value = 4800
$ws.Cells.Item(608, 16).Value = 267

# --- Append new rows 609-610 (the pair that rolled off the bottom of the window) ---
# Row 609
$ws.Cells.Item(609, 1).Value = 1
$ws.Cells.Item(609, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(609, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(609, 4).Value = 44312
$ws.Cells.Item(609, 5).Value = 15
$ws.Cells.Item(609, 6).Value = 100112033
$ws.Cells.Item(609, 7).Value = 'Lechuga'
$ws.Cells.Item(609, 8).Value = 'Escarola'
$ws.Cells.Item(609, 9).Value = 'Primera'
$ws.Cells.Item(609, 10).Value = 120
$ws.Cells.Item(609, 11).Value = 4500
$ws.Cells.Item(609, 12).Value = 5000
$ws.Cells.Item(609, 13).Value = 4750
$ws.Cells.Item(609, 14).Value = '$/caja 12 unidades'
$ws.Cells.Item(609, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(609, 16).Value = 396
$ws.Cells.Item(609, 17).Value = 12
$ws.Cells.Item(609, 18).Value = 'Hortaliza'
$ws.Cells.Item(609, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 610
$ws.Cells.Item(610, 1).Value = 1
$ws.Cells.Item(610, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(610, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(610, 4).Value = 44312
$ws.Cells.Item(610, 5).Value = 15
$ws.Cells.Item(610, 6).Value = 100112033
$ws.Cells.Item(610, 7).Value = 'Lechuga'
$ws.Cells.Item(610, 8).Value = 'Escarola'
$ws.Cells.Item(610, 9).Value = 'Segunda'
$ws.Cells.Item(610, 10).Value = 120
$ws.Cells.Item(610, 11).Value = 4500
$ws.Cells.Item(610, 12).Value = 5000
$ws.Cells.Item(610, 13).Value = 4750
$ws.Cells.Item(610, 14).Value = '$/caja 18 unidades'
$ws.Cells.Item(610, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(610, 16).Value = 264
$ws.Cells.Item(610, 17).Value = 18
$ws.Cells.Item(610, 18).Value = 'Hortaliza'
$ws.Cells.Item(610, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Dimension/used-range updates automatically on save ---
